$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Affixes")

# --- Row 381: refine the Vampire relic/enchantment description wording ---
$ws.Range("B381").Value2 = "This enchantment is best suited for a Vampire as it increases their stats and class bonuses as well as gives other Vampire specific enhancements. This is known as a kill enchantment. That is the damage on it will kill in conjunction with the life stealing."

# --- Numeric corrections on existing rows ---
$ws.Range("AM396").Value2 = 7000000000   # was 7000000
$ws.Range("X405").Value2 = 150000         # was 150000000
$ws.Range("D406").Value2 = 2.3            # was 2
$ws.Range("E406").Value2 = 2.3            # was 2
$ws.Range("G406").Value2 = 2.3            # was 2

# --- New kingdom affix rows 407-418 ---
# Row 407
$ws.Range("A407").Value2 = "Heretics Insanity"
$ws.Range("B407").Value2 = "This enchantment is best suited for a Heretic as it increases their stats and class bonuses as well as gives other Heretic specific enhancements"
$ws.Range("C407").Value2 = 1.2
$ws.Range("G407").Value2 = 1
$ws.Range("J407").Value2 = 1.1000000000000001
$ws.Range("K407").Value2 = 0
$ws.Range("L407").Value2 = 1.2
$ws.Range("M407").Value2 = 0
$ws.Range("N407").Value2 = 0
$ws.Range("O407").Value2 = 0
$ws.Range("P407").Value2 = 0
$ws.Range("Q407").Value2 = 0
$ws.Range("R407").Value2 = 0
$ws.Range("S407").Value2 = 0
$ws.Range("T407").Value2 = 0
$ws.Range("V407").Value2 = 0.3
$ws.Range("W407").Value2 = 1500000
$ws.Range("X407").Value2 = 170000
$ws.Range("Y407").Value2 = 320
$ws.Range("Z407").Value2 = 360
$ws.Range("AE407").Value2 = 1
$ws.Range("AF407").Value2 = 0
$ws.Range("AG407").Value2 = 0
$ws.Range("AH407").Value2 = 0.2
$ws.Range("AI407").Value2 = 0
$ws.Range("AL407").Value2 = 1
$ws.Range("AM407").Value2 = 17000000000
$ws.Range("AN407").Value2 = "prefix"
$ws.Range("AO407").Value2 = 0.15
$ws.Range("AP407").Value2 = 0
$ws.Range("AQ407").Value2 = 0
$ws.Range("AR407").Value2 = 0.05
# Row 408
$ws.Range("A408").Value2 = "Prophet's Raging Prayer"
$ws.Range("B408").Value2 = "This enchantment is best suited for a Prophet as it increases their stats and class bonuses as well as gives other Prophet specific enhancements"
$ws.Range("D408").Value2 = 1.5
$ws.Range("I408").Value2 = 2
$ws.Range("K408").Value2 = 0
$ws.Range("L408").Value2 = 2
$ws.Range("M408").Value2 = 0
$ws.Range("N408").Value2 = 0
$ws.Range("O408").Value2 = 0
$ws.Range("P408").Value2 = 0
$ws.Range("Q408").Value2 = 0
$ws.Range("R408").Value2 = 0
$ws.Range("S408").Value2 = 0
$ws.Range("T408").Value2 = 0
$ws.Range("V408").Value2 = 0.3
$ws.Range("W408").Value2 = 1500000
$ws.Range("X408").Value2 = 170000
$ws.Range("Y408").Value2 = 320
$ws.Range("Z408").Value2 = 360
$ws.Range("AE408").Value2 = 0
$ws.Range("AF408").Value2 = 1.6
$ws.Range("AG408").Value2 = 0
$ws.Range("AH408").Value2 = 0.6
$ws.Range("AI408").Value2 = 0
$ws.Range("AL408").Value2 = 1
$ws.Range("AM408").Value2 = 17000000000
$ws.Range("AN408").Value2 = "prefix"
$ws.Range("AO408").Value2 = 0
$ws.Range("AP408").Value2 = 0
$ws.Range("AQ408").Value2 = 0
$ws.Range("AR408").Value2 = 0
# Row 409
$ws.Range("A409").Value2 = "Earths Maddness"
$ws.Range("B409").Value2 = "This enchantment is best suited for a Ranger as it increases their stats and class bonuses as well as gives other Ranger specific enhancements"
$ws.Range("C409").Value2 = 1.2
$ws.Range("D409").Value2 = 1.05
$ws.Range("H409").Value2 = 1
$ws.Range("K409").Value2 = 1.8
$ws.Range("L409").Value2 = 0
$ws.Range("M409").Value2 = 0
$ws.Range("N409").Value2 = 0
$ws.Range("O409").Value2 = 0
$ws.Range("P409").Value2 = 0
$ws.Range("Q409").Value2 = 0
$ws.Range("R409").Value2 = 0
$ws.Range("S409").Value2 = 0
$ws.Range("T409").Value2 = 0
$ws.Range("V409").Value2 = 0.3
$ws.Range("W409").Value2 = 1500000
$ws.Range("X409").Value2 = 170000
$ws.Range("Y409").Value2 = 320
$ws.Range("Z409").Value2 = 360
$ws.Range("AE409").Value2 = 1.1000000000000001
$ws.Range("AF409").Value2 = 0.6
$ws.Range("AG409").Value2 = 0.1
$ws.Range("AH409").Value2 = 0.2
$ws.Range("AI409").Value2 = 0
$ws.Range("AL409").Value2 = 1
$ws.Range("AM409").Value2 = 17000000000
$ws.Range("AN409").Value2 = "prefix"
$ws.Range("AO409").Value2 = 0
$ws.Range("AP409").Value2 = 0
$ws.Range("AQ409").Value2 = 0
$ws.Range("AR409").Value2 = 0
# Row 410
$ws.Range("A410").Value2 = "Assassins Lucid Dream"
$ws.Range("B410").Value2 = "This enchantment is best suited for a Thief as it increases their stats and class bonuses as well as gives other Thief specific enhancements"
$ws.Range("C410").Value2 = 1.1000000000000001
$ws.Range("H410").Value2 = 1.1000000000000001
$ws.Range("K410").Value2 = 1.9
$ws.Range("L410").Value2 = 0
$ws.Range("M410").Value2 = 0
$ws.Range("N410").Value2 = 0
$ws.Range("O410").Value2 = 0
$ws.Range("P410").Value2 = 0
$ws.Range("Q410").Value2 = 0
$ws.Range("R410").Value2 = 0
$ws.Range("S410").Value2 = 0
$ws.Range("T410").Value2 = 0
$ws.Range("V410").Value2 = 0.3
$ws.Range("W410").Value2 = 1500000
$ws.Range("X410").Value2 = 170000
$ws.Range("Y410").Value2 = 320
$ws.Range("Z410").Value2 = 360
$ws.Range("AA410").Value2 = "Looting"
$ws.Range("AC410").Value2 = 1
$ws.Range("AD410").Value2 = 0.5
$ws.Range("AE410").Value2 = 1
$ws.Range("AF410").Value2 = 0
$ws.Range("AG410").Value2 = 0
$ws.Range("AH410").Value2 = 0.2
$ws.Range("AI410").Value2 = 0
$ws.Range("AL410").Value2 = 1
$ws.Range("AM410").Value2 = 17000000000
$ws.Range("AN410").Value2 = "prefix"
$ws.Range("AO410").Value2 = 0
$ws.Range("AP410").Value2 = 0
$ws.Range("AQ410").Value2 = 0
$ws.Range("AR410").Value2 = 0
# Row 411
$ws.Range("A411").Value2 = "Vampiric Ascension"
$ws.Range("B411").Value2 = "This enchantment is best suited for a Vampire as it increases their stats and class bonuses as well as gives other Vampire specific enhancements. This is known as a kill enchantment. That is the damage on it will kill in conjunction with the life stealing."
$ws.Range("C411").Value2 = 0.8
$ws.Range("D411").Value2 = 0.6
$ws.Range("G411").Value2 = 1.05
$ws.Range("K411").Value2 = 0
$ws.Range("L411").Value2 = 0
$ws.Range("M411").Value2 = 0
$ws.Range("N411").Value2 = 0
$ws.Range("O411").Value2 = 0
$ws.Range("P411").Value2 = 0
$ws.Range("Q411").Value2 = 0
$ws.Range("R411").Value2 = 0
$ws.Range("S411").Value2 = 0
$ws.Range("T411").Value2 = 0
$ws.Range("U411").Value2 = 0.7
$ws.Range("V411").Value2 = 0.3
$ws.Range("W411").Value2 = 1500000
$ws.Range("X411").Value2 = 170000
$ws.Range("Y411").Value2 = 320
$ws.Range("Z411").Value2 = 360
$ws.Range("AE411").Value2 = 0.2
$ws.Range("AF411").Value2 = 0
$ws.Range("AG411").Value2 = 0
$ws.Range("AH411").Value2 = 0.2
$ws.Range("AI411").Value2 = 0
$ws.Range("AL411").Value2 = 1
$ws.Range("AM411").Value2 = 17000000000
$ws.Range("AN411").Value2 = "prefix"
$ws.Range("AO411").Value2 = 0
$ws.Range("AP411").Value2 = 0
$ws.Range("AQ411").Value2 = 0
$ws.Range("AR411").Value2 = 0
# Row 412
$ws.Range("A412").Value2 = "Prophets Deliverance"
$ws.Range("B412").Value2 = "This enchantment is best suited for a Prophet as it increases their stats and class bonuses as well as gives other Prophet specific enhancements"
$ws.Range("D412").Value2 = 1.2
$ws.Range("I412").Value2 = 1.1000000000000001
$ws.Range("K412").Value2 = 0
$ws.Range("L412").Value2 = 1
$ws.Range("M412").Value2 = 0
$ws.Range("N412").Value2 = 0
$ws.Range("O412").Value2 = 0
$ws.Range("P412").Value2 = 0
$ws.Range("Q412").Value2 = 0
$ws.Range("R412").Value2 = 0
$ws.Range("S412").Value2 = 0
$ws.Range("T412").Value2 = 0
$ws.Range("V412").Value2 = 0.3
$ws.Range("W412").Value2 = 1500000
$ws.Range("X412").Value2 = 170000
$ws.Range("Y412").Value2 = 320
$ws.Range("Z412").Value2 = 360
$ws.Range("AE412").Value2 = 0
$ws.Range("AF412").Value2 = 1.6
$ws.Range("AG412").Value2 = 0
$ws.Range("AH412").Value2 = 0.2
$ws.Range("AI412").Value2 = 0
$ws.Range("AL412").Value2 = 1
$ws.Range("AM412").Value2 = 17000000000
$ws.Range("AN412").Value2 = "prefix"
$ws.Range("AO412").Value2 = 0
$ws.Range("AP412").Value2 = 0
$ws.Range("AQ412").Value2 = 0
$ws.Range("AR412").Value2 = 0
# Row 413
$ws.Range("A413").Value2 = "Soldiers Strike of Faith"
$ws.Range("B413").Value2 = "This enchantment is best suited for a Fighter as it increases their stats and class bonuses as well as gives other Fighter specific enhancements."
$ws.Range("C413").Value2 = 1.2
$ws.Range("E413").Value2 = 1.8
$ws.Range("F413").Value2 = 2.8
$ws.Range("G413").Value2 = 1.5
$ws.Range("H413").Value2 = 1.1000000000000001
$ws.Range("K413").Value2 = 0
$ws.Range("L413").Value2 = 0
$ws.Range("M413").Value2 = 1.1000000000000001
$ws.Range("N413").Value2 = 1.1000000000000001
$ws.Range("O413").Value2 = 1.1000000000000001
$ws.Range("P413").Value2 = 0
$ws.Range("Q413").Value2 = 0
$ws.Range("R413").Value2 = 0
$ws.Range("S413").Value2 = 0
$ws.Range("T413").Value2 = 1
$ws.Range("V413").Value2 = 0
$ws.Range("W413").Value2 = 2000000
$ws.Range("X413").Value2 = 180000
$ws.Range("Y413").Value2 = 330
$ws.Range("Z413").Value2 = 360
$ws.Range("AE413").Value2 = 1
$ws.Range("AF413").Value2 = 0
$ws.Range("AG413").Value2 = 0.55000000000000004
$ws.Range("AH413").Value2 = 0.23
$ws.Range("AI413").Value2 = 0
$ws.Range("AL413").Value2 = 1
$ws.Range("AM413").Value2 = 18000000000
$ws.Range("AN413").Value2 = "prefix"
$ws.Range("AO413").Value2 = 0
$ws.Range("AP413").Value2 = 0.1
$ws.Range("AQ413").Value2 = 0
$ws.Range("AR413").Value2 = 0.15
# Row 414
$ws.Range("A414").Value2 = "Fanatics Rage"
$ws.Range("B414").Value2 = "This enchantment is best suited for a Heretic as it increases their stats and class bonuses as well as gives other Heretic specific enhancements."
$ws.Range("C414").Value2 = 1.2
$ws.Range("J414").Value2 = 1.1000000000000001
$ws.Range("K414").Value2 = 0
$ws.Range("L414").Value2 = 1.1000000000000001
$ws.Range("M414").Value2 = 0
$ws.Range("N414").Value2 = 0
$ws.Range("O414").Value2 = 0
$ws.Range("P414").Value2 = 0
$ws.Range("Q414").Value2 = 1.1000000000000001
$ws.Range("R414").Value2 = 0
$ws.Range("S414").Value2 = 1
$ws.Range("T414").Value2 = 1
$ws.Range("V414").Value2 = 0.5
$ws.Range("W414").Value2 = 2000000
$ws.Range("X414").Value2 = 180000
$ws.Range("Y414").Value2 = 330
$ws.Range("Z414").Value2 = 360
$ws.Range("AE414").Value2 = 1
$ws.Range("AF414").Value2 = 0
$ws.Range("AG414").Value2 = 0
$ws.Range("AH414").Value2 = 0.23
$ws.Range("AI414").Value2 = 0
$ws.Range("AL414").Value2 = 1
$ws.Range("AM414").Value2 = 18000000000
$ws.Range("AN414").Value2 = "prefix"
$ws.Range("AO414").Value2 = 0.3
$ws.Range("AP414").Value2 = 0.1
$ws.Range("AQ414").Value2 = 0
$ws.Range("AR414").Value2 = 0.1
# Row 415
$ws.Range("A415").Value2 = "Clerics Heavenly Devotion"
$ws.Range("B415").Value2 = "This enchantment is best suited for a Prophet as it increases their stats and class bonuses as well as gives other Prophet specific enhancements."
$ws.Range("D415").Value2 = 1
$ws.Range("I415").Value2 = 1.1000000000000001
$ws.Range("K415").Value2 = 0
$ws.Range("L415").Value2 = 1.1000000000000001
$ws.Range("M415").Value2 = 0
$ws.Range("N415").Value2 = 0
$ws.Range("O415").Value2 = 0
$ws.Range("P415").Value2 = 1
$ws.Range("Q415").Value2 = 0
$ws.Range("R415").Value2 = 0
$ws.Range("S415").Value2 = 0.95
$ws.Range("T415").Value2 = 1
$ws.Range("V415").Value2 = 0.5
$ws.Range("W415").Value2 = 2000000
$ws.Range("X415").Value2 = 180000
$ws.Range("Y415").Value2 = 330
$ws.Range("Z415").Value2 = 360
$ws.Range("AE415").Value2 = 0
$ws.Range("AF415").Value2 = 0.65
$ws.Range("AG415").Value2 = 0
$ws.Range("AH415").Value2 = 0.23
$ws.Range("AI415").Value2 = 0
$ws.Range("AL415").Value2 = 1
$ws.Range("AM415").Value2 = 18000000000
$ws.Range("AN415").Value2 = "prefix"
$ws.Range("AO415").Value2 = 0
$ws.Range("AP415").Value2 = 0.1
$ws.Range("AQ415").Value2 = 0
$ws.Range("AR415").Value2 = 0.1
# Row 416
$ws.Range("A416").Value2 = "Natures Seeking Shadows"
$ws.Range("B416").Value2 = "This enchantment is best suited for a Heretic as it increases their stats and class bonuses as well as gives other Heretic specific enhancements."
$ws.Range("C416").Value2 = 1.1000000000000001
$ws.Range("D416").Value2 = 1
$ws.Range("H416").Value2 = 1.1000000000000001
$ws.Range("K416").Value2 = 0.8
$ws.Range("L416").Value2 = 0
$ws.Range("M416").Value2 = 0
$ws.Range("N416").Value2 = 0
$ws.Range("O416").Value2 = 1.1000000000000001
$ws.Range("P416").Value2 = 0
$ws.Range("Q416").Value2 = 0
$ws.Range("R416").Value2 = 0.95
$ws.Range("S416").Value2 = 0
$ws.Range("T416").Value2 = 1
$ws.Range("V416").Value2 = 0.5
$ws.Range("W416").Value2 = 2000000
$ws.Range("X416").Value2 = 180000
$ws.Range("Y416").Value2 = 330
$ws.Range("Z416").Value2 = 360
$ws.Range("AE416").Value2 = 1.3
$ws.Range("AF416").Value2 = 0.6
$ws.Range("AG416").Value2 = 0
$ws.Range("AH416").Value2 = 0.23
$ws.Range("AI416").Value2 = 0.28000000000000003
$ws.Range("AL416").Value2 = 1
$ws.Range("AM416").Value2 = 18000000000
$ws.Range("AN416").Value2 = "prefix"
$ws.Range("AO416").Value2 = 0
$ws.Range("AP416").Value2 = 0.1
$ws.Range("AQ416").Value2 = 0
$ws.Range("AR416").Value2 = 0.1
# Row 417
$ws.Range("A417").Value2 = "Thieves Courage"
$ws.Range("B417").Value2 = "This enchantment is best suited for a Thief as it increases their stats and class bonuses as well as gives other Thief specific enhancements."
$ws.Range("C417").Value2 = 1.1000000000000001
$ws.Range("H417").Value2 = 1.1000000000000001
$ws.Range("K417").Value2 = 1
$ws.Range("L417").Value2 = 0
$ws.Range("M417").Value2 = 0
$ws.Range("N417").Value2 = 0
$ws.Range("O417").Value2 = 1
$ws.Range("P417").Value2 = 0
$ws.Range("Q417").Value2 = 0
$ws.Range("R417").Value2 = 1.1000000000000001
$ws.Range("S417").Value2 = 0
$ws.Range("T417").Value2 = 1
$ws.Range("V417").Value2 = 0.5
$ws.Range("W417").Value2 = 2000000
$ws.Range("X417").Value2 = 180000
$ws.Range("Y417").Value2 = 330
$ws.Range("Z417").Value2 = 360
$ws.Range("AA417").Value2 = "Looting"
$ws.Range("AC417").Value2 = 1
$ws.Range("AD417").Value2 = 0.75
$ws.Range("AE417").Value2 = 1.3
$ws.Range("AF417").Value2 = 0
$ws.Range("AG417").Value2 = 0
$ws.Range("AH417").Value2 = 0.28000000000000003
$ws.Range("AI417").Value2 = 0.32
$ws.Range("AL417").Value2 = 1
$ws.Range("AM417").Value2 = 18000000000
$ws.Range("AN417").Value2 = "prefix"
$ws.Range("AO417").Value2 = 0.1
$ws.Range("AP417").Value2 = 0.1
$ws.Range("AQ417").Value2 = 0.35
$ws.Range("AR417").Value2 = 0.15
# Row 418
$ws.Range("A418").Value2 = "Vampires Imortal Blood Lust"
$ws.Range("B418").Value2 = "This enchantment is best suited for a Vampire as it increases their stats and class bonuses as well as gives other Vampire specific enhancements. This is considered a kill enchantment. That is, with the damage from life stealing and the damage from the enchantment its self, is enough to kill (in most cases)"
$ws.Range("C418").Value2 = 1.2
$ws.Range("D418").Value2 = 1
$ws.Range("G418").Value2 = 2.8
$ws.Range("K418").Value2 = 0
$ws.Range("L418").Value2 = 0
$ws.Range("M418").Value2 = 0
$ws.Range("N418").Value2 = 2.1
$ws.Range("O418").Value2 = 0
$ws.Range("P418").Value2 = 0
$ws.Range("Q418").Value2 = 0
$ws.Range("R418").Value2 = 0
$ws.Range("S418").Value2 = 0
$ws.Range("T418").Value2 = 1
$ws.Range("U418").Value2 = 0.95
$ws.Range("V418").Value2 = 0
$ws.Range("W418").Value2 = 2000000
$ws.Range("X418").Value2 = 180000
$ws.Range("Y418").Value2 = 330
$ws.Range("Z418").Value2 = 360
$ws.Range("AE418").Value2 = 1.3
$ws.Range("AF418").Value2 = 0.65
$ws.Range("AG418").Value2 = 0
$ws.Range("AH418").Value2 = 0.28000000000000003
$ws.Range("AI418").Value2 = 0
$ws.Range("AL418").Value2 = 1
$ws.Range("AM418").Value2 = 18000000000
$ws.Range("AN418").Value2 = "prefix"
$ws.Range("AO418").Value2 = 0
$ws.Range("AP418").Value2 = 0.1
$ws.Range("AQ418").Value2 = 0
$ws.Range("AR418").Value2 = 0.1

# --- Restore cursor/selection state to match the reviewed cell ---
$ws.Range("A382").Select()
